$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the culture_collection column (Z) entirely - header cell, shared string,
# and all cells to its right shift one column to the left.
$ws.Columns("Z").Delete()

# Comments are NOT auto-shifted by the column delete above, so re-home each
# note at columns Z (26) onward to the text that originally belonged to the
# next column over (mirrors the cell-data shift).
$ws.Range("Z15").Comment.Text('date of most recent douche')
$ws.Range("AA15").Comment.Text('ethnicity of the subject')
$ws.Range("AB15").Comment.Text('Plasmids that have significance phenotypic consequence')
$ws.Range("AC15").Comment.Text('history of gynecological disorders; can include multiple disorders')
$ws.Range("AD15").Comment.Text('Health or disease status of sample at time of collection')
$ws.Range("AE15").Comment.Text('Age of host at the time of sampling')
$ws.Range("AF15").Comment.Text('body mass index of the host, calculated as weight/(height)squared')
$ws.Range("AG15").Comment.Text('substance produced by the host, e.g. stool, mucus, where the sample was obtained from')
$ws.Range("AH15").Comment.Text('core body temperature of the host when sample was collected')
$ws.Range("AI15").Comment.Text('type of diet depending on the sample for animals omnivore, herbivore etc., for humans high-fat, meditteranean etc.; can include multiple diet types')
$ws.Range("AJ15").Comment.Text('Name of relevant disease, e.g. Salmonella gastroenteritis. For the controlled vocabulary, please see Human Disease Ontology, http://bioportal.bioontology.org/ontologies/1009 or MeSH, http://www.ncbi.nlm.nih.gov/mesh')
$ws.Range("AK15").Comment.Delete()
$ws.Range("AM15").AddComment('the height of subject')
$ws.Range("AN15").Comment.Text('content of last meal and time since feeding; can include multiple values')
$ws.Range("AO15").Comment.Text('most frequent job performed by subject')
$ws.Range("AP15").Comment.Delete()
$ws.Range("AQ15").AddComment('resting pulse of the host, measured as beats per minute')
$ws.Range("AR15").Comment.Text('Gender or physical sex of the host')
$ws.Range("AS15").Comment.Text('a unique identifier by which each subject can be referred to, de-identified, e.g. #131')
$ws.Range("AT15").Comment.Text('NCBI taxonomy ID of the host, e.g. 9606')
$ws.Range("AU15").Comment.Text('Type of tissue the initial sample was taken from. Controlled vocabulary, http://bioportal.bioontology.org/ontologies/1005')
$ws.Range("AV15").Comment.Text('total mass of the host at collection, the unit depends on host')
$ws.Range("AW15").Comment.Text('whether subject had hormone replacement theraphy, and if yes start date')
$ws.Range("AX15").Comment.Text('specification of whether hysterectomy was performed')
$ws.Range("AY15").Comment.Text('can include multiple medication codes')
$ws.Range("AZ15").Comment.Text('Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.')
$ws.Range("BA15").Comment.Text('A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html')
$ws.Range("BB15").Comment.Text('whether full medical history was collected')
$ws.Range("BC15").Comment.Text('date of most recent menstruation')
$ws.Range("BD15").Comment.Text('date of onset of menopause')
$ws.Range("BE15").Comment.Text('any other measurement performed or parameter collected, that is not listed here')
$ws.Range("BF15").Comment.Text('total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts')
$ws.Range("BG15").Comment.Text('oxygenation status of sample')
$ws.Range("BH15").Comment.Text('To what is the entity pathogenic')
$ws.Range("BI15").Comment.Text('type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types')
$ws.Range("BJ15").Comment.Text('date due of pregnancy')
$ws.Range("BK15").Comment.Text('Primary publication or genome report in the form of pubmed ID, DOI or URL')
$ws.Range("BL15").Comment.Text('Method or device employed for collecting sample')
$ws.Range("BM15").Comment.Text('Processing applied to the sample during or after isolation')
$ws.Range("BN15").Comment.Text('salinity of sample, i.e. measure of total salt concentration')
$ws.Range("BO15").Comment.Text('Amount or size of sample (volume, mass or area) that was collected')
$ws.Range("BP15").Comment.Text('duration for which sample was stored')
$ws.Range("BQ15").Comment.Text('location at which sample was stored, usually name of a specific freezer/room')
$ws.Range("BR15").Comment.Text('temperature at which sample was stored, e.g. -80')
$ws.Range("BS15").Comment.Text('volume (mL) or weight (g) of sample processed for DNA extraction')
$ws.Range("BT15").Comment.Text('current sexual partner and frequency of sex')
$ws.Range("BU15").Comment.Text('unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.')
$ws.Range("BV15").Comment.Text('Identifier for the physical specimen. Use format: "[<institution-code>:[<collection-code>:]]<specimen_id>", eg, "UAM:Mamm:52179". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a ''structured voucher''. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier')
$ws.Range("BW15").Comment.Text('Information about the genetic distinctness of the lineage (eg., biovar, serovar)')
$ws.Range("BX15").Comment.Text('temperature of the sample at time of sampling')
$ws.Range("BY15").Comment.Text('Feeding position in food chain (eg., chemolithotroph)')
$ws.Range("BZ15").Comment.Text('history of urogenital disorders, can include multiple disorders')
$ws.Range("CA15").Comment.Delete()
